# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Sandia" (Vega Modelo de Temuco) right
# before the previous last rows of data, shifting the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 404 - this pushes the existing
# rows 404:422 down to 408:426 and grows the used range to A1:R426.
$ws.Rows("404:407").Insert()

# Shared, constant columns for every data row in this sheet.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112028
$categoria = "Sandia"
$variedad  = "Sin especificar"
$kgOUnid   = 1
$clasif    = "Hortaliza"

function Set-Row($r, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidadComercializacion, $origen, $precioKg) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $catId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $variedad
    $ws.Cells.Item($r, 9).Value  = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 14).Value = $unidadComercializacion
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $precioKg
    $ws.Cells.Item($r, 17).Value = $kgOUnid
    $ws.Cells.Item($r, 18).Value = $clasif
}

Set-Row 404 44578 "Extra"   1500 3000 3000 3000 "$/unidad" "Región del Maule" 3000
Set-Row 405 44578 "Primera" 5000 2500 2500 2500 "$/unidad" "Región del Maule" 2500
Set-Row 406 44578 "Segunda" 4000 2000 2000 2000 "$/unidad" "Región del Maule" 2000
Set-Row 407 44578 "Tercera" 1500 1000 1000 1000 "$/unidad" "Región del Maule" 1000
